$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G (header "K") values regenerated to use K instead of Strike#
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 4
$ws.Range("G4").Value = 3
$ws.Range("G5").Value = 4
$ws.Range("G6").Value = 3
$ws.Range("G7").Value = 2
